# Applies the crypto price / 1h-volume refresh described in the commit
# ("Updated cryptos list ... with GitHub Actions"), plus the NEARProtocol /
# ApeXProtocol row swap (rows 47-48).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
# The "Price" column holds numeric-looking values (e.g. "104.80", "51.704.26")
# as literal text. Force a Text number format on each cell before writing it,
# otherwise Excel auto-converts it to a number and silently drops the trailing
# zero / collapses the thousands-dot formatting.
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '51.704.26'
$ws.Range('E2').Value = '  +1.35%  '

# Row 3
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.991.30'
$ws.Range('E3').Value = '  +2.77%  '

# Row 4
$ws.Range('E4').Value = '  +0.08%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '385.95'
$ws.Range('E5').Value = '  +3.10%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '104.80'
$ws.Range('E6').Value = '  +3.23%  '

# Row 7
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.547'
$ws.Range('E7').Value = '  +0.74%  '

# Row 8
$ws.Range('E8').Value = '  +0.13%  '

# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.597'
$ws.Range('E9').Value = '  +2.21%  '

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '37.34'
$ws.Range('E10').Value = '  +1.56%  '

# Row 11
$ws.Range('E11').Value = '  +0.23%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.0849'
$ws.Range('E12').Value = '  +1.86%  '

# Row 13
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '3.466.99'
$ws.Range('E13').Value = '  +2.89%  '

# Row 14
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '18.43'
$ws.Range('E14').Value = '  +1.10%  '

# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.63'
$ws.Range('E15').Value = '  +3.72%  '

# Row 16
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '2.993.44'
$ws.Range('E16').Value = '  +2.96%  '

# Row 17
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '1.01'
$ws.Range('E17').Value = '  +9.66%  '

# Row 18
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '51.688.97'
$ws.Range('E18').Value = '  +1.39%  '

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '3.29'
$ws.Range('E19').Value = '  +1.52%  '

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '7.47'
$ws.Range('E20').Value = '  +3.82%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '12.95'
$ws.Range('E21').Value = '  +0.98%  '

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.0₃0968'
$ws.Range('E22').Value = '  +2.82%  '

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '69.19'
$ws.Range('E23').Value = '  +1.50%  '

# Row 24
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '263.74'
$ws.Range('E24').Value = '  +1.70%  '

# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.93'
$ws.Range('E25').Value = '  +8.97%  '

# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '8.37'
$ws.Range('E26').Value = '  +18.41%  '

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '7.69'
$ws.Range('E27').Value = '  +21.65%  '

# Row 28
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.171'
$ws.Range('E28').Value = '  +1.89%  '

# Row 29
$ws.Range('E29').Value = '  +14.31%  '

# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '26.10'
$ws.Range('E30').Value = '  +1.87%  '

# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.999'
$ws.Range('E31').Value = '  -0.14%  '

# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '9.92'
$ws.Range('E32').Value = '  +0.95%  '

# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '35.06'
$ws.Range('E33').Value = '  +2.96%  '

# Row 34
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '51.07'
$ws.Range('E34').Value = '  -0.22%  '

# Row 35
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.07'
$ws.Range('E35').Value = '  -1.75%  '

# Row 36
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.0455'
$ws.Range('E36').Value = '  +7.74%  '

# Row 37
$ws.Range('E37').Value = '  -0.13%  '

# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '3.04'
$ws.Range('E38').Value = '  +2.10%  '

# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '17.11'
$ws.Range('E39').Value = '  +1.00%  '

# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.61'
$ws.Range('E40').Value = '  +1.21%  '

# Row 41
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.85'
$ws.Range('E41').Value = '  +0.96%  '

# Row 42
$ws.Range('E42').Value = '  +3.59%  '

# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '122.52'
$ws.Range('E43').Value = '  +2.58%  '

# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '21.90'
$ws.Range('E44').Value = '  +0.36%  '

# Row 45
$ws.Range('E45').Value = '  +18.14%  '

# Row 46
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.05'
$ws.Range('E46').Value = '  -1.63%  '

# Row 47
$ws.Range('B47').Value = 'ApeXProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.37'
$ws.Range('E47').Value = '  +2.75%  '

# Row 48
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '3.30'
$ws.Range('E48').Value = '  +5.40%  '

# Row 49
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.041.61'
$ws.Range('E49').Value = '  +1.44%  '

# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0335'
$ws.Range('E50').Value = '  +8.94%  '

# Row 51
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.865'
$ws.Range('E51').Value = '  +2.45%  '
